# Complete restructure and rewrite of documentation ready for v2
# -----------------------------------------------------------------
$wb = $excel.ActiveWorkbook

$wsNotes   = $wb.Worksheets.Item("Notes")
$wsStudies = $wb.Worksheets.Item("studies")
$wsSurveys = $wb.Worksheets.Item("surveys")
$wsCounts  = $wb.Worksheets.Item("counts")

# -----------------------------------------------------------------
# 1. studies sheet - rebuild header / data row (keep F2 hyperlink cell
#    untouched since value+style+hyperlink target are unchanged there)
# -----------------------------------------------------------------
$wsStudies.Range("A1:E2").Clear()
$wsStudies.Range("G1:G2").Clear()

$wsStudies.Range("A1").Value = "study_id"
$wsStudies.Range("B1").Value = "study_label"
$wsStudies.Range("C1").Value = "description"
$wsStudies.Range("D1").Value = "access_level"
$wsStudies.Range("E1").Value = "contributors"
$wsStudies.Range("F1").Value = "reference"
$wsStudies.Range("G1").Value = "reference_year"

$wsStudies.Range("A2").Value = "foo"
$wsStudies.Range("D2").Value = "public"

# -----------------------------------------------------------------
# 2. surveys sheet - rebuild header / data row, inserting two new
#    columns (location_method / time_method) amongst the existing set
# -----------------------------------------------------------------
$wsSurveys.Range("A1:M2").Clear()

$wsSurveys.Range("A1").Value = "study_id"
$wsSurveys.Range("B1").Value = "survey_id"
$wsSurveys.Range("C1").Value = "country_name"
$wsSurveys.Range("D1").Value = "site_name"
$wsSurveys.Range("E1").Value = "latitude"
$wsSurveys.Range("F1").Value = "longitude"
$wsSurveys.Range("G1").Value = "location_method"
$wsSurveys.Range("H1").Value = "location_notes"
$wsSurveys.Range("I1").Value = "collection_start"
$wsSurveys.Range("J1").Value = "collection_end"
$wsSurveys.Range("K1").Value = "collection_day"
$wsSurveys.Range("L1").Value = "time_method"
$wsSurveys.Range("M1").Value = "time_notes"

$wsSurveys.Range("A2").Value = "foo"
$wsSurveys.Range("B2").Value = "S01"
$wsSurveys.Range("E2").Value = 0
$wsSurveys.Range("F2").Value = 0
$wsSurveys.Range("H2").Value = "example data"
$wsSurveys.Range("K2").Value = "2020-01-01"
$wsSurveys.Range("M2").Value = "example data"

# Columns I:K already inherit the "text" style (index 2) from the
# column-level style definition. L needs the same "text" style applied
# explicitly (copy format from a cell that already carries it).
$wsSurveys.Range("I1").Copy() | Out-Null
$wsSurveys.Range("L1").PasteSpecial(-4122) | Out-Null
$wsSurveys.Range("K2").Copy() | Out-Null
$wsSurveys.Range("L2").PasteSpecial(-4122) | Out-Null
$wsSurveys.Range("L2").Value = ""
$excel.CutCopyMode = 0

# -----------------------------------------------------------------
# 3. counts sheet - same table layout, renamed key columns / value
# -----------------------------------------------------------------
$wsCounts.Range("A1").Value = "study_id"
$wsCounts.Range("B1").Value = "survey_id"

$wsCounts.Range("A2").Value = "foo"
$wsCounts.Range("A3").Value = "foo"

# -----------------------------------------------------------------
# 4. Selections / active cells per sheet
# -----------------------------------------------------------------
$wsNotes.Range("D15").Select() | Out-Null
$wsStudies.Range("D3").Select() | Out-Null
$wsSurveys.Range("A1:M2").Select() | Out-Null
$wsCounts.Range("B2").Select() | Out-Null

# -----------------------------------------------------------------
# 5. Active sheet -> studies (2nd tab)
# -----------------------------------------------------------------
$wsStudies.Activate()
